$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.631.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -8.41%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.649.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -9.39%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.73%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'219.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -6.02%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.5104"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -13.57%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.62%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.2527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -7.73%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'21.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.17%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.06109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -10.01%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07369"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.88%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.635.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -10.51%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'4.465"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.32%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.5717"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -8.28%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'1.871.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -9.47%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.000008032"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -14.06%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'64.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -13.86%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'26.628.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -7.63%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'4.960"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -8.33%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'1.015"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.99%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'10.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -7.23%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'182.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -12.34%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'1.013"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.55%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'6.167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -8.90%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'142.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -7.61%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'7.570"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.60%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  -10.24%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'14.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -7.80%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'1.326"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.88%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.05725"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -10.43%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.333"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -6.78%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'3.424"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -7.55%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  -7.34%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.572"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.03%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.9793"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.57%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'2.436"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.79%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'0.5926"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.11%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'2.616"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.05%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.8639"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.90%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01563"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -8.22%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'1.064.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.77%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.016"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.06%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'5.701"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -11.70%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'95.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.97%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'1.775.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -10.15%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00000000109"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.34%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'1.015"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.37%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'55.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -8.10%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4384"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.09%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.05206"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.11%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'7.814"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.84%  "
$ws.Range("E51").Style = "Normal"
